# Generate Report for Handoff
# Updates the localization-status workbook to reflect a fresh handoff run:
# the "Ready for handoff" rows get a new Priority ("ht" instead of "low")
# and their handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (col G) for the
# "Ready for handoff" rows (4-7) moves to the new generation time.
foreach ($r in 4..7) {
    $overview.Cells.Item($r, 7).Value = "2016-08-28 04:31:28"
}

# zh-cn sheet: Priority (col E) low -> ht, Latest Handoff Datetime (col H)
# refreshed for the same rows.
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-28 04:31:24"
}

# de-de sheet: Priority (col E) low -> ht.
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
}
